# Generate Report for Archive
#
# 1. The localization status "Ready for handoff" moved on to "In Translation"
#    for the single tracked file. This value shows up on the Overview sheet
#    (columns E/F, one per locale) as well as on each locale sheet's
#    "Status" column (column C).
# 2. The now-shorter status text let the "Status" columns shrink, so their
#    column widths are narrowed to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.43
$overview.Columns.Item(6).ColumnWidth = 12.43

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.43

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.43
